$wb = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison": update Seasonality Index (L) and MyForecast (D) values ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("L2").Value = 1.07
$ws1.Range("L3").Value = 1.08
$ws1.Range("L4").Value = 0.89
$ws1.Range("L5").Value = 0.91
$ws1.Range("L6").Value = 0.86
$ws1.Range("L7").Value = 0.95
$ws1.Range("L8").Value = 0.96

$ws1.Range("D9").Value = 8
$ws1.Range("L9").Value = 0.98

$ws1.Range("D10").Value = 8
$ws1.Range("L10").Value = 0.97

$ws1.Range("L11").Value = 0.99

$ws1.Range("D12").Value = 8
$ws1.Range("L12").Value = 0.86

$ws1.Range("D13").Value = 8
$ws1.Range("L13").Value = 1.09

$ws1.Range("L14").Value = 0.85
$ws1.Range("L15").Value = 0.89
$ws1.Range("L16").Value = 0.93
$ws1.Range("L17").Value = 1.07

# --- Sheet "Summary": update Total Forecast (16 Weeks) value, keep it stored as text ---
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "139"
